$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "99.003.71"
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("D3").Value = "3.302.61"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'254.63"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.40%  "
$ws.Range("D6").Value = "'624.19"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("E7").Value = "  +29.34%  "
$ws.Range("D8").Value = "'0.408"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +5.87%  "
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("D10").Value = "'0.969"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +21.50%  "
$ws.Range("D11").Value = "3.299.27"
$ws.Range("E11").Value = "  -1.30%  "
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("D13").Value = "'39.86"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +11.96%  "
$ws.Range("D14").Value = "98.655.92"
$ws.Range("E14").Value = "  +1.30%  "
$ws.Range("E15").Value = "  +1.37%  "
$ws.Range("D16").Value = "3.922.69"
$ws.Range("E16").Value = "  -1.01%  "
$ws.Range("D17").Value = "'5.47"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.07%  "
$ws.Range("D18").Value = "3.305.15"
$ws.Range("E18").Value = "  -1.16%  "
$ws.Range("D19").Value = "'3.46"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.43%  "
$ws.Range("D20").Value = "'15.45"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.22%  "
$ws.Range("D21").Value = "'6.29"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +8.21%  "
$ws.Range("D22").Value = "'485.09"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.46%  "
$ws.Range("D23").Value = "'9.39"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.07%  "
$ws.Range("D24").Value = "'0.0000202"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.62%  "
$ws.Range("D25").Value = "'5.62"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").Value = "'88.93"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("D27").Value = "'11.93"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.29%  "
$ws.Range("E28").Value = "  +29.92%  "
$ws.Range("D29").Value = "3.487.81"
$ws.Range("E29").Value = "  -1.07%  "
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("E31").Value = "  +11.65%  "
$ws.Range("E32").Value = "  +2.48%  "
$ws.Range("D33").Value = "'10.21"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +10.39%  "
$ws.Range("D34").Value = "'0.994"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.63%  "
$ws.Range("D35").Value = "'27.78"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.45%  "
$ws.Range("D36").Value = "'0.470"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +4.66%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.148"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.96%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").Value = "'7.19"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.80%  "
$ws.Range("D39").Value = "'1.94"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("D40").Value = "'24.81"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("D41").Value = "'488.73"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.74%  "
$ws.Range("D42").Value = "'3.62"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.62%  "
$ws.Range("E43").Value = "  -3.57%  "
$ws.Range("D44").Value = "'0.783"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.61%  "
$ws.Range("D46").Value = "'3.11"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -6.31%  "
$ws.Range("E47").Value = "  +1.90%  "
$ws.Range("D48").Value = "'158.39"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.64%  "
$ws.Range("D49").Value = "'7.32"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +15.66%  "
$ws.Range("D50").Value = "'0.846"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.22%  "
$ws.Range("D51").Value = "'4.71"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +4.42%  "
